$d = $word.ActiveDocument

# 1. Remove the stray "V " before the city placeholder in the signature line.
$d.Content.Find.Execute(" V {{ city }}, dne ", $true, $false, $false, $false, $false, $true, 1, $false, " {{ city }}, dne ", 2)

# 2. Nudge the floating "Text Box 2" signature box: Word re-flowed it a hair
#    after the text edit above, so the anchor's vertical offset and the
#    shape's stored extents shift by a few hundred EMU, and the VML
#    fallback's style string picks up the same (rounded) values.
$p = $d.Paragraphs.Item(16)
$r = $p.Range
$xml = $r.WordOpenXML
$xml = $xml.Replace('<wp:posOffset>174625</wp:posOffset>', '<wp:posOffset>182880</wp:posOffset>')
$xml = $xml.Replace('<wp:extent cx="2661920" cy="715645"/>', '<wp:extent cx="2662555" cy="715645"/>')
$xml = $xml.Replace('<a:ext cx="2661120" cy="714960"/>', '<a:ext cx="2661840" cy="714960"/>')
$xml = $xml.Replace('margin-left:290.65pt;margin-top:13.75pt;width:209.5pt;height:56.25pt', 'margin-left:290.65pt;margin-top:14.4pt;width:209.55pt;height:56.25pt')
$r.InsertXML($xml) | Out-Null
